$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 73: new CRM sample row, 16th sample "from OR" ---
# Copy the date-format style from A72 down to A73 so the new date cell
# matches the existing column-A date formatting (reuses style index 1
# instead of creating a new numFmt entry).
$ws.Range("A72").Copy() | Out-Null
$ws.Range("A73").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A73").Value = 43792
$ws.Range("B73").Value = 2195.1166635711702
$ws.Range("C73").Value = 2207.0300000000002
$ws.Range("D73").Formula = "=100*(B73-C73)/C73"
$ws.Range("E73").Value = 169
$ws.Range("F73").Value = "Crm opened 11/19/2020"

# --- Rows 74-76: trailing formula-only rows (blank inputs -> #DIV/0!) ---
$ws.Range("D74").Formula = "=100*(B74-C74)/C74"
$ws.Range("D75").Formula = "=100*(B75-C75)/C75"
$ws.Range("D76").Formula = "=100*(B76-C76)/C76"

# --- View state: scroll position + active selection on F72:F73 ---
$ws.Range("F72:F73").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 52
$excel.ActiveWindow.ScrollColumn = 3
